# "Add files via upload" — refresh the 8.2.1 indicator metadata sheet:
#   * the organisation's website address is corrected
#   * the indicator title text is tweaked ("одного" -> "каждого", no colon)
#   * the active cell selection is left where the author's edit ended up

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 ("Сайт организации (если есть)"): www.stat.kg -> www.stat.gov.kg.
# Nudge the font so the edited cell gets its own (un-wrapped) style record,
# same as happens when a cell is retyped by hand in Excel.
$ws.Range("B10").Value = "www.stat.gov.kg"
$ws.Range("B10").Font.Name = "Calibri"

# Row 4 ("Индикатор"): drop the colon and swap "одного" -> "каждого".
$ws.Range("B4").Value = "8.2.1 Ежегодные темпы роста реального ВВП на каждого занятого"

# Leave the selection on B8, matching the saved workbook view.
$ws.Range("B8").Select()
